$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix duplicated "类" typos in header row 1
$ws.Range("B1").Value = "其他食品类居民消费价格指数(上年=100)"
$ws.Range("D1").Value = "奶类居民消费价格指数(上年=100)"
$ws.Range("E1").Value = "干鲜瓜果类居民消费价格指数(上年=100)"
$ws.Range("H1").Value = "畜肉类居民消费价格指数(上年=100)"
$ws.Range("I1").Value = "禽肉类居民消费价格指数(上年=100)"
$ws.Range("K1").Value = "糖果糕点类居民消费价格指数(上年=100)"
$ws.Range("N1").Value = "薯类居民消费价格指数(上年=100)"
$ws.Range("O1").Value = "蛋类居民消费价格指数(上年=100)"
$ws.Range("Q1").Value = "豆类居民消费价格指数(上年=100)"

# Add row 7: 2021年 (copy formatting of the preceding year's label cell)
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = "2021年"
$ws.Range("B7").Value = 100.4
$ws.Range("C7").Value = 102
$ws.Range("D7").Value = 101.8
$ws.Range("E7").Value = 102.1
$ws.Range("F7").Value = 109.4
$ws.Range("G7").Value = 101.5
$ws.Range("H7").Value = 82.8
$ws.Range("I7").Value = 96.8
$ws.Range("J7").Value = 101.1
$ws.Range("K7").Value = 101.4
$ws.Range("L7").Value = 101.1
$ws.Range("M7").Value = 105
$ws.Range("N7").Value = 99.7
$ws.Range("O7").Value = 110.8
$ws.Range("P7").Value = 101.2
$ws.Range("Q7").Value = 106.6
$ws.Range("R7").Value = 99.7
$ws.Range("S7").Value = 98.59999999999999
$ws.Range("T7").Value = 106.9
$ws.Range("U7").Value = 102.8
$ws.Range("V7").Value = 105.6

# Add row 8: 2022年, only R8 has a value
$ws.Range("A6").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "2022年"
$ws.Range("R8").Value = 102.4
